# Automated update of mapa_interactivo_Optical_Power.xlsx
# - Remove the completed "BLANCO ENCALADA 5402" case (old row 34)
# - Add a new pending case "ALBARIÑO 1331" (new row 17), shifting the
#   rows in between down by one
# - Update the "BLANCO ENCALADA 5618" case (row 61) to the new
#   "VALLESE, FELIPE 490" decommission case
# - Remove the two completed "BLANCO ENCALADA 5432" / "BLANCO ENCALADA 5188"
#   cases (old rows 62-63)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Remove the old row 34 (its data is obsolete / resolved) -------------
$ws.Range("A34").EntireRow.Delete()

# --- Insert a brand-new row at 17 for the new case ------------------------
$ws.Range("A17").EntireRow.Insert()

$ws.Cells.Item(17, 1).Value  = "'6066"
$ws.Cells.Item(17, 2).Value  = "'2/26/2025"
$ws.Cells.Item(17, 3).Value  = "ALBARIÑO 1331"
$ws.Cells.Item(17, 4).Value  = "'9"
$ws.Cells.Item(17, 5).Value  = "'803651213"
$ws.Cells.Item(17, 6).Value  = "Optical Power"
$ws.Cells.Item(17, 7).Value  = "Pendiente"
$ws.Cells.Item(17, 8).Value  = "Podrida en la base"
$ws.Cells.Item(17, 9).Value  = "'1"
$ws.Cells.Item(17, 10).Value = "Cambio"
$ws.Cells.Item(17, 11).Value = "Sin equipos"
$ws.Cells.Item(17, 12).Value = "Pasante"
$ws.Cells.Item(17, 13).Value = -58.496255
$ws.Cells.Item(17, 14).Value = -34.650599

# --- Update row 61 with the new decommission case -------------------------
$ws.Cells.Item(61, 1).Value  = "'6130"
$ws.Cells.Item(61, 3).Value  = "VALLESE, FELIPE 490"
$ws.Cells.Item(61, 4).Value  = "'6"
$ws.Cells.Item(61, 5).Value  = "'"
$ws.Cells.Item(61, 8).Value  = "Desmontar columa ya traspasaron un nodo solo falta retirar"
$ws.Cells.Item(61, 10).Value = "Desmonte"
$ws.Cells.Item(61, 13).Value = -58.440448
$ws.Cells.Item(61, 14).Value = -34.611223

# --- Remove the two trailing, now-resolved rows 62 and 63 -----------------
$ws.Range("A62:A63").EntireRow.Delete()
